$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates to the "Price" (D) and "Volume(1h)" (E) columns,
# reflecting refreshed cryptocurrency market data.
$updates = @(
    @{ Row = 2; D = '67.713.38'; E = '  -0.52%  ' },
    @{ Row = 3; D = '3.259.58'; E = '  -0.26%  ' },
    @{ Row = 4; D = $null; E = '  -0.03%  ' },
    @{ Row = 5; D = '580.74'; E = '  -0.93%  ' },
    @{ Row = 6; D = '185.02'; E = '  +0.28%  ' },
    @{ Row = 8; D = '0.607'; E = '  +1.38%  ' },
    @{ Row = 9; D = '3.253.54'; E = '  -0.39%  ' },
    @{ Row = 10; D = $null; E = '  -2.21%  ' },
    @{ Row = 11; D = '6.58'; E = '  -2.06%  ' },
    @{ Row = 12; D = $null; E = '  -0.84%  ' },
    @{ Row = 13; D = '3.816.34'; E = '  -0.55%  ' },
    @{ Row = 14; D = $null; E = '  -0.06%  ' },
    @{ Row = 15; D = '27.72'; E = '  -3.04%  ' },
    @{ Row = 16; D = '67.700.22'; E = '  -0.52%  ' },
    @{ Row = 17; D = $null; E = '  -1.22%  ' },
    @{ Row = 18; D = '3.264.87'; E = '  -0.12%  ' },
    @{ Row = 19; D = $null; E = '  -1.76%  ' },
    @{ Row = 20; D = '13.60'; E = '  +0.17%  ' },
    @{ Row = 21; D = '394.61'; E = '  +3.31%  ' },
    @{ Row = 22; D = '7.62'; E = '  -1.49%  ' },
    @{ Row = 23; D = '71.52'; E = $null },
    @{ Row = 24; D = $null; E = '  +0.10%  ' },
    @{ Row = 25; D = $null; E = '  +0.71%  ' },
    @{ Row = 26; D = $null; E = '  -1.85%  ' },
    @{ Row = 27; D = $null; E = '  -2.43%  ' },
    @{ Row = 28; D = '9.60'; E = '  -1.39%  ' },
    @{ Row = 29; D = $null; E = '  +0.71%  ' },
    @{ Row = 30; D = $null; E = '  -2.01%  ' },
    @{ Row = 31; D = '5.55'; E = '  -4.32%  ' },
    @{ Row = 32; D = '22.70'; E = '  -0.96%  ' },
    @{ Row = 33; D = $null; E = '  -2.38%  ' },
    @{ Row = 34; D = $null; E = '  -2.26%  ' },
    @{ Row = 35; D = $null; E = '  +0.05%  ' },
    @{ Row = 36; D = '162.01'; E = '  -0.80%  ' },
    @{ Row = 37; D = $null; E = '  -3.88%  ' },
    @{ Row = 38; D = '1.91'; E = '  +2.04%  ' },
    @{ Row = 39; D = '26.63'; E = '  +0.09%  ' },
    @{ Row = 40; D = $null; E = '  -3.47%  ' },
    @{ Row = 41; D = $null; E = '  -1.21%  ' },
    @{ Row = 42; D = '6.48'; E = '  -4.04%  ' },
    @{ Row = 43; D = '2.49'; E = '  -5.19%  ' },
    @{ Row = 44; D = $null; E = '  -0.03%  ' },
    @{ Row = 45; D = '40.68'; E = $null },
    @{ Row = 46; D = '2.613.84'; E = '  -0.51%  ' },
    @{ Row = 47; D = '24.83'; E = '  -2.66%  ' },
    @{ Row = 48; D = '335.12'; E = '  -1.57%  ' },
    @{ Row = 49; D = $null; E = '  -2.04%  ' },
    @{ Row = 50; D = '6.38'; E = '  +2.02%  ' },
    @{ Row = 51; D = $null; E = '  -0.40%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D = Price
        $origStyle = $cell.Style
        # Force text so values such as "580.74" or "13.60" are not
        # reinterpreted as numbers; this keeps the cell type as a string,
        # matching the original inline-string cells.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = $origStyle
    }

    if ($null -ne $u.E) {
        $eCell = $ws.Cells.Item($row, 5)  # column E = Volume(1h)
        $eCell.Value = $u.E
    }
}
